$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 11111343
$ws.Range("I6").Value = 20833416
$ws.Range("K6").Value = 62500248
$ws.Range("M6").Value = -62500136
# Row 32
$ws.Range("H32").Value = 7070
$ws.Range("I32").Value = 7133.6665
$ws.Range("J32").Value = 6974.5
$ws.Range("K32").Value = 7133.6665
$ws.Range("L32").Value = 6974.5
$ws.Range("M32").Value = -6807.6665
$ws.Range("N32").Value = -7626.5
# Row 63
$ws.Range("H63").Value = 47271
$ws.Range("J63").Value = 47271
$ws.Range("L63").Value = 47271
$ws.Range("N63").Value = -48519
# Row 66
$ws.Range("H66").Value = 47271
$ws.Range("J66").Value = 47271
$ws.Range("L66").Value = 141813
$ws.Range("N66").Value = -148053
# Row 103
$ws.Range("H103").Value = 38099.91
$ws.Range("I103").Value = 200202
$ws.Range("K103").Value = 600606
$ws.Range("M103").Value = -600020

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4632.3335
$ws.Range("I61").Value = 4948.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 4948.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -4736.5
$ws.Range("N61").Value = -4424
# Row 63
$ws.Range("H63").Value = 2833
$ws.Range("J63").Value = 2833
$ws.Range("L63").Value = 2833
$ws.Range("N63").Value = -4205
# Row 66
$ws.Range("H66").Value = 2833
$ws.Range("J66").Value = 2833
$ws.Range("L66").Value = 14165
$ws.Range("N66").Value = -21029
# Row 74
$ws.Range("H74").Value = 2407.5
$ws.Range("I74").Value = 1865.625
$ws.Range("J74").Value = 4575
$ws.Range("K74").Value = 1865.625
$ws.Range("L74").Value = 4575
$ws.Range("M74").Value = -991.625
$ws.Range("N74").Value = -6323
# Row 77
$ws.Range("H77").Value = 2407.5
$ws.Range("I77").Value = 1865.625
$ws.Range("J77").Value = 4575
$ws.Range("K77").Value = 9328.125
$ws.Range("L77").Value = 22875
$ws.Range("M77").Value = -4960.125
$ws.Range("N77").Value = -31611
# Row 132
$ws.Range("H132").Value = 2650.8235
$ws.Range("I132").Value = 2185.4482
$ws.Range("K132").Value = 6556.344599999999
$ws.Range("M132").Value = -4026.344599999999
# Row 136
$ws.Range("H136").Value = 4632.3335
$ws.Range("I136").Value = 4948.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 14845.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -12295.5
$ws.Range("N136").Value = -17100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6508.2856
$ws.Range("I86").Value = 10414.857
$ws.Range("J86").Value = 2601.7144
$ws.Range("K86").Value = 10414.857
$ws.Range("L86").Value = 2601.7144
$ws.Range("M86").Value = -9291.857
$ws.Range("N86").Value = -4847.7144
# Row 89
$ws.Range("H89").Value = 6508.2856
$ws.Range("I89").Value = 10414.857
$ws.Range("J89").Value = 2601.7144
$ws.Range("K89").Value = 52074.285
$ws.Range("L89").Value = 13008.572
$ws.Range("M89").Value = -46458.285
$ws.Range("N89").Value = -24240.572
# Row 105
$ws.Range("H105").Value = 2425.1177
$ws.Range("I105").Value = 2103.0715
$ws.Range("K105").Value = 2103.0715
$ws.Range("M105").Value = -356.0715

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7990.8
$ws.Range("I31").Value = 7486.2
$ws.Range("K31").Value = 7486.2
$ws.Range("M31").Value = -7191.2
# Row 34
$ws.Range("H34").Value = 7990.8
$ws.Range("I34").Value = 7486.2
$ws.Range("K34").Value = 7486.2
$ws.Range("M34").Value = -7284.2
# Row 59
$ws.Range("H59").Value = 24499.5
# Row 60
$ws.Range("H60").Value = 17000
$ws.Range("I60").Value = 15000
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -14489
$ws.Range("N60").Value = -21022
# Row 62
$ws.Range("H62").Value = 21210.715
$ws.Range("I62").Value = 21695
$ws.Range("K62").Value = 21695
$ws.Range("M62").Value = -21071
# Row 65
$ws.Range("H65").Value = 21210.715
$ws.Range("I65").Value = 21695
$ws.Range("K65").Value = 108475
$ws.Range("M65").Value = -105355
# Row 74
$ws.Range("H74").Value = 58999.332
$ws.Range("J74").Value = 59249.5
$ws.Range("L74").Value = 59249.5
$ws.Range("N74").Value = -60997.5
# Row 77
$ws.Range("H77").Value = 58999.332
$ws.Range("J77").Value = 59249.5
$ws.Range("L77").Value = 177748.5
$ws.Range("N77").Value = -186484.5
# Row 86
$ws.Range("H86").Value = 11413.934
$ws.Range("J86").Value = 13199
$ws.Range("L86").Value = 13199
$ws.Range("N86").Value = -15445
# Row 89
$ws.Range("H89").Value = 11413.934
$ws.Range("J89").Value = 13199
$ws.Range("L89").Value = 65995
$ws.Range("N89").Value = -77227
# Row 132
$ws.Range("H132").Value = 27977.4
$ws.Range("I132").Value = 8397.263000000001
$ws.Range("K132").Value = 25191.789
$ws.Range("M132").Value = -22661.789

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").Value = 0
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").Value = 0
# Row 131
$ws.Range("H131").Value = 12501757
$ws.Range("J131").Value = 1777.7106
$ws.Range("L131").Value = 5333.1318
$ws.Range("N131").Value = -15413.1318
# Row 132
$ws.Range("H132").Value = 85923.164
$ws.Range("J132").Value = 128447.25
$ws.Range("L132").Value = 1156025.25
$ws.Range("N132").Value = -1161085.25
# Row 140
$ws.Range("H140").Value = 2635.0435
$ws.Range("I140").Value = 2635.0435
$ws.Range("K140").Value = 7905.130500000001
$ws.Range("M140").Value = -2725.130500000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 13506000
$ws.Range("I18").Value = 27000000
$ws.Range("K18").Value = 27000000
$ws.Range("M18").Value = -26999707
# Row 35
$ws.Range("H35").Value = 14971.667
$ws.Range("I35").Value = 12915
$ws.Range("K35").Value = 12915
$ws.Range("M35").Value = -12617
# Row 80
$ws.Range("H80").Value = 2766.0908
$ws.Range("I80").Value = 2505.8572
$ws.Range("K80").Value = 2505.8572
$ws.Range("M80").Value = -1507.8572
# Row 83
$ws.Range("H83").Value = 2766.0908
$ws.Range("I83").Value = 2505.8572
$ws.Range("K83").Value = 12529.286
$ws.Range("M83").Value = -7537.286
# Row 126
$ws.Range("H126").Value = 30433.666
$ws.Range("I126").Value = 43983.25
$ws.Range("J126").Value = 19594
$ws.Range("K126").Value = 131949.75
$ws.Range("L126").Value = 58782
$ws.Range("M126").Value = -129479.75
$ws.Range("N126").Value = -63722

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
# Row 55
$ws.Range("H55").Value = 875.4828
$ws.Range("I55").Value = 830.1905
$ws.Range("J55").Value = 994.375
$ws.Range("K55").Value = 830.1905
$ws.Range("L55").Value = 994.375
$ws.Range("M55").Value = -657.1905
$ws.Range("N55").Value = -1340.375
# Row 82
$ws.Range("H82").Value = 2021.6207
$ws.Range("I82").Value = 2401.8125
$ws.Range("J82").Value = 1553.6923
$ws.Range("K82").Value = 2401.8125
$ws.Range("L82").Value = 1553.6923
$ws.Range("M82").Value = -2040.8125
$ws.Range("N82").Value = -2275.6923
# Row 85
$ws.Range("H85").Value = 2021.6207
$ws.Range("I85").Value = 2401.8125
$ws.Range("J85").Value = 1553.6923
$ws.Range("K85").Value = 2401.8125
$ws.Range("L85").Value = 1553.6923
$ws.Range("M85").Value = -1153.8125
$ws.Range("N85").Value = -4049.6923
# Row 100
$ws.Range("H100").Value = 7608.5454
$ws.Range("I100").Value = 11298.333
$ws.Range("K100").Value = 11298.333
$ws.Range("M100").Value = -10757.333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2723.5
$ws.Range("I107").Value = 3855.6
$ws.Range("J107").Value = 2094.5557
$ws.Range("K107").Value = 11566.8
$ws.Range("L107").Value = 6283.6671
$ws.Range("M107").Value = -9646.799999999999
$ws.Range("N107").Value = -10123.6671
# Row 113
$ws.Range("H113").Value = 8399.700000000001
$ws.Range("I113").Value = 5999.3076
$ws.Range("K113").Value = 17997.9228
$ws.Range("M113").Value = -15827.9228
# Row 132
$ws.Range("H132").Value = 8975.547
$ws.Range("I132").Value = 10668.361
$ws.Range("K132").Value = 32005.083
$ws.Range("M132").Value = -29475.083
